$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Hydrogen" / "Iron & steel" results cell for Steel
$ws.Range("B3").Value = 99275.06305729281

# Float re-serialization nudge to match the committed value for D8
$ws.Range("D8").Value = 40707.5009778319
